# This workbook keeps a rolling weekly log of price observations for
# "Betarraga" (beet) at Feria Lagunitas de Puerto Montt. Each weekly
# refresh inserts a brand-new record at row 103 (the top of the dated
# series) and pushes every subsequent row down by one, with the record
# that used to be the very last one (row 221) surviving onto a newly
# appended row 222.
#
# Columns D (Fecha) and J:Q (Volumen .. Kg o Unidades) are the ones that
# vary row to row; A,B,C,E,F,G,H,I,R stay constant across the whole
# block, so shifting the full A:R block down by one row (then overwriting
# row 103 with the new observation) reproduces the expected result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 103
$lastRow  = 221
$newLastRow = $lastRow + 1

# 1) Capture the complete original rows firstRow..lastRow (all columns,
#    A:R) before making any changes, so we still have the untouched data
#    once we start overwriting it.
$srcRange = $ws.Range("A$($firstRow):R$($lastRow)")
$srcValues = $srcRange.Value2

# 2) Shift everything down by one row: row (firstRow+1) gets what used to
#    be in row firstRow, ..., row newLastRow (222) gets what used to be in
#    row lastRow (221).
$dstRange = $ws.Range("A$($firstRow + 1):R$($newLastRow)")
$dstRange.Value2 = $srcValues

# Row newLastRow needs the Fecha column's date number format (it was a
# brand-new row, so it currently has no special formatting yet).
$ws.Range("D$newLastRow").NumberFormat = $ws.Range("D$lastRow").NumberFormat

# 3) Overwrite row 103 with the brand-new observation for this week.
$ws.Range("D$firstRow").Value2 = 44539
$ws.Range("J$firstRow").Value2 = 500
$ws.Range("K$firstRow").Value2 = 900
$ws.Range("L$firstRow").Value2 = 1000
$ws.Range("M$firstRow").Value2 = 950
$ws.Range("N$firstRow").Value2 = "$/paquete 5 unidades"
$ws.Range("O$firstRow").Value2 = "Región del Maule"
$ws.Range("P$firstRow").Value2 = 190
$ws.Range("Q$firstRow").Value2 = 5

$ws.Range("A1").Select() | Out-Null
